$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text labels (shared strings) in row 9
$ws.Range("F9").Value = "Tasa de seroreversion"
$ws.Range("G9").Value = "Tasa de serereversión Rhat"

# Update the active selection to C12
$ws.Range("C12").Select()

# Update column F width (closest achievable value to the target 19.6640625)
$ws.Columns.Item(6).ColumnWidth = 18.83
